$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update header H1: TEnrolled -> TAG ---
$ws.Range("H1").Value = "TAG"

# --- Update existing rows 2-4: column H becomes a short text tag instead of a numeric enrollment count ---
$ws.Range("H2").Value = "o"
$ws.Range("H3").Value = "o"
$ws.Range("H4").Value = "o"

# --- Build new rows 5-8 by copying the formatting of existing rows, then overwrite values ---
# Row 5 (Chess / ONGOING=1) -> formatted like row 2 (hyperlink style in column C)
$ws.Range("A2:H2").Copy($ws.Range("A5:H5"))
# Row 6 (Valorant / ONGOING=1) -> formatted like row 2 as well
$ws.Range("A2:H2").Copy($ws.Range("A6:H6"))
# Row 7 (BGMI / COMPLETED=1) -> formatted like row 4 (plain text style in column C, not hyperlink style)
$ws.Range("A4:H4").Copy($ws.Range("A7:H7"))
# Row 8 (Carrom / COMPLETED=1) -> formatted like row 2 again
$ws.Range("A2:H2").Copy($ws.Range("A8:H8"))

# Row 5: Chess, UPCOMING=0, ONGOING=1, COMPLETED=0, TAG=i
$ws.Range("A5").Value = "Chess"
$ws.Range("B5").Value = "Chess"
$ws.Range("C5").Value = "https://cdn.wallpapersafari.com/3/7/SFe72B.jpg"
$ws.Range("D5").Value = "2023-10-17T11:25:00"
$ws.Range("E5").Value = 0
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0
$ws.Range("H5").Value = "i"

# Row 6: Valorant, UPCOMING=0, ONGOING=1, COMPLETED=0, TAG=e
$ws.Range("A6").Value = "Valorant"
$ws.Range("B6").Value = "Valorant"
$ws.Range("C6").Value = "https://hdqwalls.com/wallpapers/basketball-hd.jpg"
$ws.Range("D6").Value = "2023-10-17T11:25:00"
$ws.Range("E6").Value = 0
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0
$ws.Range("H6").Value = "e"

# Row 7: BGMI, UPCOMING=0, ONGOING=0, COMPLETED=1, TAG=e
$ws.Range("A7").Value = "BGMI"
$ws.Range("B7").Value = "BGMI"
$ws.Range("C7").Value = "https://wallpaperdig.com/wp-content/uploads/2021/02/XDCFVBGNHJM.jpg"
$ws.Range("D7").Value = "2023-10-17T11:25:00"
$ws.Range("E7").Value = 0
$ws.Range("F7").Value = 0
$ws.Range("G7").Value = 1
$ws.Range("H7").Value = "e"

# Row 8: Carrom, UPCOMING=0, ONGOING=0, COMPLETED=1, TAG=i
$ws.Range("A8").Value = "Carrom"
$ws.Range("B8").Value = "Carrom"
$ws.Range("C8").Value = "https://cdn.wallpapersafari.com/3/7/SFe72B.jpg"
$ws.Range("D8").Value = "2023-10-17T11:25:00"
$ws.Range("E8").Value = 0
$ws.Range("F8").Value = 0
$ws.Range("G8").Value = 1
$ws.Range("H8").Value = "i"

# G8 carries no explicit cell style in the target workbook - strip formatting on it.
$ws.Range("G8").ClearFormats()

# --- Add the two new hyperlinks, then restore the worksheet's custom hyperlink style ---
$ws.Hyperlinks.Add($ws.Range("C5"), "https://cdn.wallpapersafari.com/3/7/SFe72B.jpg")
$ws.Range("C2").Copy($ws.Range("C5"))
$ws.Range("C5").Value = "https://cdn.wallpapersafari.com/3/7/SFe72B.jpg"

$ws.Hyperlinks.Add($ws.Range("C8"), "https://cdn.wallpapersafari.com/3/7/SFe72B.jpg")
$ws.Range("C2").Copy($ws.Range("C8"))
$ws.Range("C8").Value = "https://cdn.wallpapersafari.com/3/7/SFe72B.jpg"

# --- Two trailing placeholder rows (9 & 10), empty but carrying the column C / D formatting ---
$ws.Range("C2").Copy($ws.Range("C9"))
$ws.Range("C9").ClearContents()
$ws.Range("D2").Copy($ws.Range("C10"))
$ws.Range("C10").ClearContents()

# --- Final selection matches the authored workbook ---
$ws.Range("G8").Select()
